$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F ("想去人数") that must be
# updated identically on both the "展览" sheet (sheet1) and the
# "全部类型" sheet (sheet4).
$updates = @{
    2  = 180
    4  = 12494
    6  = 146
    7  = 30
    8  = 89
    10 = 197
    11 = 451
    17 = 4074
    18 = 97
    19 = 8
    20 = 941
    23 = 57
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
